# Re-generate the file list and review-date index.
# - A handful of rows in the 'File' list were reordered (re-generated listing order).
# - Every 'Review date' that fell on the 16th of the month is moved to the 19th,
#   except for the two rows that needed to absorb the newly-inserted 2023-01-18 entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: rows whose 'File' text moved as part of the re-generated order ---
$fileMoves = @{
    82 = 'Infection_and_sepsis/SARI/Influenza A and B Virology Sampling and Oseltamivir Dose.pdf'
    83 = 'GI_Liver_and_Transplant/Nasogastric feeding protocol.pdf'
    89 = 'GI_Liver_and_Transplant/Nasojejunal feeding protocol.pdf'
    90 = 'GI_Liver_and_Transplant/Jejunostomy feeding protocol.pdf'
    91 = 'ECLS/RIE ECLS Anti Xa Protocol.pdf'
    92 = 'Drugs/calcium.pdf'
    100 = 'End_of_life_care/CMO & NRS Guidance for Doctors completing MCCD - Sept 22.pdf'
    101 = 'Organ_donation/Organ Retrieval SOP.pdf'
    102 = 'Drugs/vasopressin_sepsis.pdf'
    103 = 'Neurological/Management of traumatic brain injury.pdf'
    105 = 'Transfer/ACCP Transfers.pdf'
    106 = 'Drugs/ketamine for status epilepticus.pdf'
    107 = 'Drugs/thiopentone.pdf'
    108 = 'Covid-19/videos/Donning and Doffing Video.pdf'
    109 = 'Drugs/valproate.pdf'
}
foreach ($row in $fileMoves.Keys) {
    $ws.Cells.Item($row, 1).Value = $fileMoves[$row]
}

# --- Column B: review dates. Re-write every date cell as literal text (never let
#     Excel reinterpret the YYYY-MM-DD text as a real date/number). ---
$lastRow = $ws.UsedRange.Rows.Count
for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 2)
    $old = $cell.Text
    if ([string]::IsNullOrEmpty($old)) { continue }

    if ($row -eq 82) {
        $new = '2023-01-18'
    } elseif ($row -eq 83) {
        $new = '2023-01-19'
    } elseif ($old.EndsWith('-16')) {
        $new = $old.Substring(0, $old.Length - 2) + '19'
    } else {
        $new = $old
    }

    if ($new -ne $old) {
        $cell.NumberFormat = '@'
        $cell.Value = $new
    }
}

